$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "revision" column (E) only ever held the literal 0 and is being
# dropped entirely; everything to its right (old column F,
# "Revised_Requirement") shifts one column to the left to become the new
# column E. Deleting the whole column (rather than just clearing cells)
# removes the now-unused "revision" shared string and shifts F -> E.
$ws.Range("E1").EntireColumn.Delete()

# Give the (now) Requirement / Revised_Requirement columns explicit,
# best-fit-style widths like the saved workbook has.
$ws.Columns.Item(4).ColumnWidth = 236.1796875
$ws.Columns.Item(5).ColumnWidth = 255.6328125

# Restore the view: 85% zoom and a selection sitting a couple of rows
# below the data, as captured in the saved workbook.
$win = $excel.ActiveWindow
$win.Zoom = 85
$ws.Range("F13").Select()
